$d = $word.ActiveDocument

$replacements = @(
    @("499×3=1497", "602×5=3010"),
    @("613×3=1839", "122×7=854"),
    @("141×8=1128", "164×9=1476"),
    @("977×8=7816", "527×2=1054"),
    @("221×6=1326", "477×2=954"),
    @("198×8=1584", "205×9=1845"),
    @("806×5=4030", "939×7=6573"),
    @("610×7=4270", "518×7=3626"),
    @("819×8=6552", "629×6=3774"),
    @("517×8=4136", "784×4=3136"),
    @("391×2=782", "953×5=4765"),
    @("881×6=5286", "833×2=1666"),
    @("538×3=1614", "650×2=1300"),
    @("692×4=2768", "381×7=2667"),
    @("314×6=1884", "798×3=2394"),
    @("473×9=4257", "547×2=1094"),
    @("152×5=760", "319×6=1914"),
    @("213×2=426", "411×8=3288"),
    @("974×6=5844", "664×3=1992"),
    @("829×3=2487", "469×6=2814"),
    @("954×2=1908", "363×6=2178"),
    @("430×5=2150", "248×4=992"),
    @("607×9=5463", "862×5=4310"),
    @("698×5=3490", "798×6=4788"),
    @("583×5=2915", "715×9=6435")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Write-Host "Done applying replacements"
